# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Tue Jan  2 15:38:42 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{ Cell = 'D2'; Value = '45.283.60'; Numeric = $false },
    @{ Cell = 'E2'; Value = '  +5.33%  '; Numeric = $false },
    @{ Cell = 'D3'; Value = '2.366.19'; Numeric = $false },
    @{ Cell = 'E3'; Value = '  +2.55%  '; Numeric = $false },
    @{ Cell = 'E4'; Value = '  -0.47%  '; Numeric = $false },
    @{ Cell = 'D5'; Value = '110.02'; Numeric = $true },
    @{ Cell = 'E5'; Value = '  +4.12%  '; Numeric = $false },
    @{ Cell = 'D6'; Value = '309.24'; Numeric = $true },
    @{ Cell = 'E6'; Value = '  -0.52%  '; Numeric = $false },
    @{ Cell = 'D7'; Value = '0.630'; Numeric = $true },
    @{ Cell = 'E7'; Value = '  +0.44%  '; Numeric = $false },
    @{ Cell = 'E8'; Value = '  -0.18%  '; Numeric = $false },
    @{ Cell = 'D9'; Value = '0.616'; Numeric = $true },
    @{ Cell = 'E9'; Value = '  +1.88%  '; Numeric = $false },
    @{ Cell = 'D10'; Value = '41.29'; Numeric = $true },
    @{ Cell = 'E10'; Value = '  +3.06%  '; Numeric = $false },
    @{ Cell = 'D11'; Value = '0.0918'; Numeric = $true },
    @{ Cell = 'E11'; Value = '  +0.65%  '; Numeric = $false },
    @{ Cell = 'D12'; Value = '8.50'; Numeric = $true },
    @{ Cell = 'E12'; Value = '  +1.56%  '; Numeric = $false },
    @{ Cell = 'E13'; Value = '  +1.38%  '; Numeric = $false },
    @{ Cell = 'D14'; Value = '0.988'; Numeric = $true },
    @{ Cell = 'E14'; Value = '  -0.29%  '; Numeric = $false },
    @{ Cell = 'D15'; Value = '2.726.12'; Numeric = $false },
    @{ Cell = 'E15'; Value = '  +2.61%  '; Numeric = $false },
    @{ Cell = 'D16'; Value = '15.44'; Numeric = $true },
    @{ Cell = 'E16'; Value = '  +0.86%  '; Numeric = $false },
    @{ Cell = 'D17'; Value = '2.369.96'; Numeric = $false },
    @{ Cell = 'E17'; Value = '  +3.20%  '; Numeric = $false },
    @{ Cell = 'D18'; Value = '45.216.18'; Numeric = $false },
    @{ Cell = 'E18'; Value = '  +5.48%  '; Numeric = $false },
    @{ Cell = 'D19'; Value = '7.32'; Numeric = $true },
    @{ Cell = 'E19'; Value = '  -2.39%  '; Numeric = $false },
    @{ Cell = 'D20'; Value = '0.0000107'; Numeric = $true },
    @{ Cell = 'E20'; Value = '  +1.48%  '; Numeric = $false },
    @{ Cell = 'D21'; Value = '13.11'; Numeric = $true },
    @{ Cell = 'E21'; Value = '  -2.41%  '; Numeric = $false },
    @{ Cell = 'D22'; Value = '73.54'; Numeric = $true },
    @{ Cell = 'E22'; Value = '  +0.02%  '; Numeric = $false },
    @{ Cell = 'D23'; Value = '3.45'; Numeric = $true },
    @{ Cell = 'E23'; Value = '  -0.72%  '; Numeric = $false },
    @{ Cell = 'D24'; Value = '260.87'; Numeric = $true },
    @{ Cell = 'E24'; Value = '  -2.26%  '; Numeric = $false },
    @{ Cell = 'D25'; Value = '2.28'; Numeric = $true },
    @{ Cell = 'E25'; Value = '  +1.81%  '; Numeric = $false },
    @{ Cell = 'E26'; Value = '  -0.41%  '; Numeric = $false },
    @{ Cell = 'D27'; Value = '7.45'; Numeric = $true },
    @{ Cell = 'E27'; Value = '  -4.17%  '; Numeric = $false },
    @{ Cell = 'D28'; Value = '11.09'; Numeric = $true },
    @{ Cell = 'E28'; Value = '  +1.28%  '; Numeric = $false },
    @{ Cell = 'E29'; Value = '  +2.48%  '; Numeric = $false },
    @{ Cell = 'D30'; Value = '22.47'; Numeric = $true },
    @{ Cell = 'E30'; Value = '  +0.83%  '; Numeric = $false },
    @{ Cell = 'D31'; Value = '0.0962'; Numeric = $true },
    @{ Cell = 'E31'; Value = '  +10.76%  '; Numeric = $false },
    @{ Cell = 'D32'; Value = '37.89'; Numeric = $true },
    @{ Cell = 'E32'; Value = '  -0.93%  '; Numeric = $false },
    @{ Cell = 'D33'; Value = '169.91'; Numeric = $true },
    @{ Cell = 'E33'; Value = '  +2.59%  '; Numeric = $false },
    @{ Cell = 'E34'; Value = '  +4.56%  '; Numeric = $false },
    @{ Cell = 'E35'; Value = '  +0.19%  '; Numeric = $false },
    @{ Cell = 'B36'; Value = 'Kaspa'; Numeric = $false },
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; Numeric = $false },
    @{ Cell = 'D36'; Value = '0.116'; Numeric = $true },
    @{ Cell = 'E36'; Value = '  +4.02%  '; Numeric = $false },
    @{ Cell = 'B37'; Value = 'RenderToken'; Numeric = $false },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; Numeric = $false },
    @{ Cell = 'D37'; Value = '4.80'; Numeric = $true },
    @{ Cell = 'E37'; Value = '  +3.70%  '; Numeric = $false },
    @{ Cell = 'D38'; Value = '3.00'; Numeric = $true },
    @{ Cell = 'E38'; Value = '  +6.98%  '; Numeric = $false },
    @{ Cell = 'D39'; Value = '3.94'; Numeric = $true },
    @{ Cell = 'E39'; Value = '  +8.73%  '; Numeric = $false },
    @{ Cell = 'D40'; Value = '0.0356'; Numeric = $true },
    @{ Cell = 'E40'; Value = '  -0.20%  '; Numeric = $false },
    @{ Cell = 'D41'; Value = '1.73'; Numeric = $true },
    @{ Cell = 'E41'; Value = '  +9.43%  '; Numeric = $false },
    @{ Cell = 'D42'; Value = '101.36'; Numeric = $true },
    @{ Cell = 'E42'; Value = '  -2.55%  '; Numeric = $false },
    @{ Cell = 'E43'; Value = '  +1.64%  '; Numeric = $false },
    @{ Cell = 'B44'; Value = 'MultiversX'; Numeric = $false },
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'; Numeric = $false },
    @{ Cell = 'D44'; Value = '69.99'; Numeric = $true },
    @{ Cell = 'E44'; Value = '  -1.36%  '; Numeric = $false },
    @{ Cell = 'B45'; Value = 'Celestia'; Numeric = $false },
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'; Numeric = $false },
    @{ Cell = 'D45'; Value = '13.09'; Numeric = $true },
    @{ Cell = 'E45'; Value = '  +6.55%  '; Numeric = $false },
    @{ Cell = 'D46'; Value = '0.999'; Numeric = $true },
    @{ Cell = 'E46'; Value = '  -0.65%  '; Numeric = $false },
    @{ Cell = 'D47'; Value = '82.68'; Numeric = $true },
    @{ Cell = 'E47'; Value = '  +8.56%  '; Numeric = $false },
    @{ Cell = 'D48'; Value = '9.42'; Numeric = $true },
    @{ Cell = 'E48'; Value = '  +6.23%  '; Numeric = $false },
    @{ Cell = 'D49'; Value = '112.82'; Numeric = $true },
    @{ Cell = 'E49'; Value = '  +1.32%  '; Numeric = $false },
    @{ Cell = 'D50'; Value = '5.54'; Numeric = $true },
    @{ Cell = 'E50'; Value = '  +6.57%  '; Numeric = $false },
    @{ Cell = 'B51'; Value = 'MinaProtocolToken'; Numeric = $false },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'; Numeric = $false },
    @{ Cell = 'D51'; Value = '1.55'; Numeric = $true },
    @{ Cell = 'E51'; Value = '  +6.75%  '; Numeric = $false }
)

foreach ($e in $edits) {
    $rng = $ws.Range($e.Cell)
    if ($e.Numeric) {
        # Force text storage so numeric-looking values (prices) are not
        # reinterpreted as numbers, matching the source inlineStr cells.
        $rng.NumberFormat = "@"
        $rng.Value = $e.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $e.Value
    }
}
